# Final commit for attempt 1
# Fill in the DTT Assessment Hour Log with the actual logged hours/work items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTT Test Hour Log")

# --- Rows 4-8 --------------------------------------------------------------
# The Subject/Date/Description text in these rows is already correct in the
# source file; only the "Amount of hours" column needs to change from a
# textual "0,5"-style value to a genuine number.
$ws.Range("B4").Value = 0.5
$ws.Range("B5").Value = 0.25
$ws.Range("B6").Value = 0.75
$ws.Range("B7").Value = 0.2
$ws.Range("B8").Value = 2.5

# --- Row 9 ---------------------------------------------------------------
$ws.Range("A9").Value = "Creation of getAll and getByID (User Story 2)"
$ws.Range("B9").Value = 1.5
$ws.Range("D9").Value = "Had some trouble figuring out how to get the results of a query with the given library and also figured out how to defend against sql injection. Sources: https://www.w3schools.com/php/php_mysql_prepared_statements.asp AND https://www.w3schools.com/php/php_mysql_select.asp"

# --- Row 10 --------------------------------------------------------------
$ws.Range("A10").Value = "Creation of create (User Story 2)"
$ws.Range("B10").Value = 1.25
$ws.Range("D10").Value = "Created the endpoint for creating a faclility and its tags"

# --- Row 11 --------------------------------------------------------------
$ws.Range("A11").Value = "Update of getAll and getByID (User Story 2)"
$ws.Range("B11").Value = 0.75
$ws.Range("D11").Value = "Updated the get endpoints to match with the requirement to also provide tags and location in the response"

# --- Rows 12 & 13 ----------------------------------------------------------
# Subjects for both rows are written first (matches the shared-string order
# of the saved workbook), then the descriptions - "deleting" before the
# longer "updating" one.
$ws.Range("A12").Value = "Creation of update (User Story 2)"
$ws.Range("A13").Value = "Creation of delete (User Story 2)"
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 0.25
$ws.Range("D13").Value = "Created the endpoint for deleting a facility and its tags"
$ws.Range("D12").Value = "Created the endpoint for updating a faclility and its tags, the database relations have been changed to CASCADE to make it easier to remove and (re)add the tags. I chose to take this approach in stead of getting the tags and filtering which ones are already connected to the facility is for simplicity sake. Had I had more time on my hands I could have attempted to try this."

# --- Row 14 --------------------------------------------------------------
$ws.Range("A14").Value = "Creation of search function (User Story 3)"
$ws.Range("B14").Value = 0.5
$ws.Range("D14").Value = "Created the endpoint for searching on certain query parameters"

# --- Row 15 --------------------------------------------------------------
$ws.Range("A15").Value = "Cleanup Code"
$ws.Range("B15").Value = 0.5
$ws.Range("D15").Value = "Cleaned up duplicate code and added comments here and there"

# --- Row 16 --------------------------------------------------------------
$ws.Range("A16").Value = "Report/Documentation"
$ws.Range("B16").Value = 1
$ws.Range("D16").Value = "Writing of API Documentation in Postman and export of SQL dump"

# Dates for the newly-filled rows (9-16 continue the date progression of
# the rest of the log).
$ws.Range("C9").Value = 45315
$ws.Range("C10").Value = 45316
$ws.Range("C11").Value = 45316
$ws.Range("C12").Value = 45316
$ws.Range("C13").Value = 45316
$ws.Range("C14").Value = 45316
$ws.Range("C15").Value = 45316
$ws.Range("C16").Value = 45316

# Recalculate the running total.
$wb.Application.Calculate()

# Update the selection to match the saved view.
$ws.Range("H17").Select()
